$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text updates ---
# "Volume 32   Number  31" -> "Volume 32   Number  32" (last run "31" -> "32")
$ws.Range("A8").Characters(21, 2).Text = "32"

# "Report Covering the Week  7/28/2025  Through  8/3/2025"
#  -> "Report Covering the Week  8/4/2025  Through  8/10/2025"
# Replace the later run first so the earlier runs character offsets stay valid.
$ws.Range("C9").Characters(47, 8).Text = "8/10/2025"
$ws.Range("C9").Characters(27, 9).Text = "8/4/2025"

# --- Weekly crime-stat table updates (rows 15-31) ---

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("N15").Value = -42.857142857142

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 89
$ws.Range("J16").Value = 81
$ws.Range("K16").Value = 9.876543209876
$ws.Range("L16").Value = -11.881188118811
$ws.Range("M16").Value = -47.337278106508
$ws.Range("N16").Value = -85.117056856187

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = -9.090909090909
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -17.073170731707
$ws.Range("I17").Value = 271
$ws.Range("J17").Value = 303
$ws.Range("K17").Value = -10.561056105610
$ws.Range("L17").Value = -4.577464788732
$ws.Range("M17").Value = 43.386243386243
$ws.Range("N17").Value = -43.892339544513

# Row 18
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 64
$ws.Range("K18").Value = -12.5
$ws.Range("L18").Value = -35.632183908046
$ws.Range("M18").Value = -66.467065868263
$ws.Range("N18").Value = -94.895168641750

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 115.384615384615
$ws.Range("I19").Value = 195
$ws.Range("J19").Value = 238
$ws.Range("K19").Value = -18.067226890756
$ws.Range("L19").Value = -22.310756972111
$ws.Range("M19").Value = -22.924901185770
$ws.Range("N19").Value = -43.641618497109

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 36.363636363636
$ws.Range("I20").Value = 43
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = -48.192771084337
$ws.Range("L20").Value = -51.685393258427
$ws.Range("M20").Value = -65.6
$ws.Range("N20").Value = -94.743276283618

# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 0
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = 15.662650602409
$ws.Range("I21").Value = 676
$ws.Range("J21").Value = 783
$ws.Range("K21").Value = -13.665389527458
$ws.Range("L21").Value = -18.847539015606
$ws.Range("M21").Value = -27.311827956989
$ws.Range("N21").Value = -80.053113012688

# Row 23
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 42
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = -16
$ws.Range("L23").Value = -30
$ws.Range("M23").Value = 61.538461538461

# Row 24
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -31.25
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = -12.676056338028
$ws.Range("I24").Value = 744
$ws.Range("J24").Value = 767
$ws.Range("K24").Value = -2.998696219035
$ws.Range("L24").Value = -0.932090545938
$ws.Range("M24").Value = -14.874141876430

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 25
$ws.Range("H25").Value = -7.407407407407
$ws.Range("I25").Value = 384
$ws.Range("J25").Value = 338
$ws.Range("K25").Value = 13.609467455621
$ws.Range("L25").Value = 31.958762886597

# Row 26
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -27.777777777777
$ws.Range("F26").Value = 60
$ws.Range("G26").Value = 71
$ws.Range("H26").Value = -15.492957746478
$ws.Range("I26").Value = 469
$ws.Range("J26").Value = 475
$ws.Range("K26").Value = -1.263157894736
$ws.Range("L26").Value = 5.630630630630
$ws.Range("M26").Value = -37.383177570093

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 300

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -7.692307692307
$ws.Range("I28").Value = 58
$ws.Range("J28").Value = 57
$ws.Range("K28").Value = 1.754385964912
$ws.Range("L28").Value = 7.407407407407

# Row 29
$ws.Range("L29").Value = -77.777777777777
$ws.Range("N29").Value = -93.846153846153

# Row 30
$ws.Range("L30").Value = -75
$ws.Range("N30").Value = -92.727272727272

# Row 31
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("C31").Value = 1
$ws.Range("I31").Value = 5
$ws.Range("K31").Value = 66.666666666666
$ws.Range("L31").Value = 400
